$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1, columns A:U): "_old" -> "_FV2310" and
#    "_new" -> "_FV2404" suffixes (column K, "diff", is left untouched).
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value2
    if ($text -like "*_old") {
        $cell.Value2 = $text.Replace("_old", "_FV2310")
    } elseif ($text -like "*_new") {
        $cell.Value2 = $text.Replace("_new", "_FV2404")
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table (ListObject) so the renamed
#    headers become the table's column headers, with filtering enabled.
# ---------------------------------------------------------------------------
$rng = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# ---------------------------------------------------------------------------
# 3) Freeze the header row (row 1) so it stays visible while scrolling.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
